# Rockwell first working version
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the Next_update value for the Schneider row (row 4) to match
# Last_update (column C) — was 45266, now 44901.
$ws.Range("D4").Value = 44901

# Move the active cell/selection to D5 (was F10).
$ws.Range("D5").Select() | Out-Null

# Reflect the updated window geometry recorded for this workbook view.
try {
    $aw = $excel.ActiveWindow
    $aw.Left = -110
    $aw.Top = -110
    $aw.Width = 38620
    $aw.Height = 21220
} catch {
    # Window geometry is host-UI state; ignore if unsupported by the runtime.
}
